$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Cell VALUES — rows 8-15 are rewritten so that rows that compare equal
#    ("c-"/"c+" matched-key rows) collapse into a single "-"/"+" pair and the
#    new row (ID 4 / Thailand) that was previously swallowed by the old
#    "=" (unchanged) row now shows up as its own "-" (old) / "+" (new) rows.
# ---------------------------------------------------------------------------

$ws.Range("A8").Value  = "-"
$ws.Range("B8").Value  = "2022-11-05 00:00:00"
$ws.Range("C8").Value  = "2"
$ws.Range("D8").Value  = "banana"
$ws.Range("E8").Value  = "Chile"
$ws.Range("F8").Value  = ""
$ws.Range("G8").Value  = "yellow"

$ws.Range("A9").Value  = "-"
$ws.Range("B9").Value  = "2022-11-20 00:00:00"
$ws.Range("C9").Value  = "3"
$ws.Range("D9").Value  = "coconut"
$ws.Range("E9").Value  = "Hawaii"
$ws.Range("G9").Value  = "brown"

$ws.Range("A10").Value = "-"
$ws.Range("B10").Value = "2022-11-21 00:00:00"
$ws.Range("C10").Value = "4"
$ws.Range("D10").Value = "dried mango"
$ws.Range("E10").Value = "Thailand"
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = "orange"

$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("G11").Value = ""

$ws.Range("A12").Value = "-"
$ws.Range("B12").Value = "Trailing row here"
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""

$ws.Range("A13").Value = "+"
$ws.Range("C13").Value = "grapes"
$ws.Range("D13").Value = "pebbles"
$ws.Range("E13").Value = "purple"

$ws.Range("A14").Value = "+"
$ws.Range("B14").Value = "2"
$ws.Range("C14").Value = "mango"
$ws.Range("D14").Value = "oval"
$ws.Range("E14").Value = "yellow"

$ws.Range("A15").Value = "+"
$ws.Range("B15").Value = "4"
$ws.Range("C15").Value = "dried mango"
$ws.Range("D15").Value = "flat"
$ws.Range("E15").Value = "orange"
$ws.Range("F15").Value = ""
$ws.Range("G15").Value = ""

# ---------------------------------------------------------------------------
# 2) Cell STYLES — reuse the existing (unchanged) rows 1/2/5 as format
#    templates so the same shared fills/cellXfs get reused instead of Excel
#    minting new ones.
# ---------------------------------------------------------------------------

$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1,B1,C1,D1,E1,F1,G1,A7,C7,D7,G7,A8,C8,D8,G8,A9,C9,D9,G9,A10,C10,D10,G10,A11,B11,C11,D11,E11,F11,G11,A12,B12,C12,D12,E12,F12,G12").PasteSpecial(-4122)

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A2,B2,C2,D2,E2,F2,G2,A13,B13,C13,D13,E13,F13,G13,A14,B14,C14,D14,E14,F14,G14,A15,B15,C15,D15,E15,F15,G15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Recolor the shared "-"/"+" row fills (the diff-tool's old/new palette)
#    and the "c-"/"c+" accent fills. Every cell that carries a given fill is
#    repainted together so the color change applies once per palette slot.
# ---------------------------------------------------------------------------

# old-row pink fill: 00FFBBBB -> 00FFB6C1
$ws.Range("A1,B1,C1,D1,E1,F1,G1,A7,C7,D7,G7,A8,C8,D8,G8,A9,C9,D9,G9,A10,C10,D10,G10,A11,B11,C11,D11,E11,F11,G11,A12,B12,C12,D12,E12,F12,G12").Interior.Color = 12695295

# new-row green fill: 00BBFFBB -> 00B6FFC1
$ws.Range("A2,B2,C2,D2,E2,F2,G2,A13,B13,C13,D13,E13,F13,G13,A14,B14,C14,D14,E14,F14,G14,A15,B15,C15,D15,E15,F15,G15").Interior.Color = 12713910

# old-row header/accent fill: 00FFDDDD -> 00FFDDE2
$ws.Range("B5,E5,B6,E6,B7,E7,B8,E8,B9,E9,B10,E10").Interior.Color = 14867967

# new-row header/accent fill: 00DDFFDD -> 00DDFFE2
$ws.Range("F5,F6,F7,F8,F9,F10").Interior.Color = 14876637

Write-Output "done"
